$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose literal text must be forced to Text type (they look like numbers/dates
# and Excel would otherwise auto-convert them to a numeric/date value).
$textCells = @{
    "H2" = "0.9992"
    "I2" = "0.5883"
    "Q2" = "2025-04-04"
    "H3" = "0.9976"
    "I3" = "0.4645"
    "Q3" = "2025-04-16"
    "H4" = "0.9952"
    "I4" = "0.3847"
    "H5" = "0.9911"
    "I5" = "0.3107"
    "Q5" = "2025-04-02"
    "H6" = "0.9843"
    "I6" = "0.2447"
    "H7" = "0.6628"
    "Q7" = "2025-06-18"
    "H8" = "0.4271"
    "Q8" = "2025-06-10"
    "H9" = "0.3690"
    "I9" = "0.9982"
    "H10" = "0.5522"
    "H11" = "0.2569"
    "I11" = "0.7320"
    "H12" = "0.2016"
    "I12" = "0.1700"
    "Q12" = "2025-06-18"
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}

# Cells that are already non-numeric-looking text; plain assignment keeps them as text.
$ws.Range("O2").Value = "R$ -375.97"
$ws.Range("P2").Value = "R$ -17.93"
$ws.Range("O3").Value = "R$ -330.90"
$ws.Range("P3").Value = "R$ -71.70"
$ws.Range("O4").Value = "R$ -366.46"
$ws.Range("P4").Value = "R$ -77.41"
$ws.Range("J5").Value = "ITM → ATM"
$ws.Range("O5").Value = "R$ -363.35"
$ws.Range("P5").Value = "R$ -69.38"
$ws.Range("O6").Value = "R$ -372.29"
$ws.Range("P6").Value = "R$ -108.90"
$ws.Range("O7").Value = "R$ 1051.01"
$ws.Range("P7").Value = "R$ 1051.01"
$ws.Range("O8").Value = "R$ 437.61"
$ws.Range("P8").Value = "R$ 960.34"
$ws.Range("O9").Value = "R$ 308.47"
$ws.Range("P9").Value = "R$ 968.36"
$ws.Range("O10").Value = "R$ 2047.65"
$ws.Range("P10").Value = "R$ 2047.65"
$ws.Range("O11").Value = "R$ 590.57"
$ws.Range("P11").Value = "R$ 724.28"
$ws.Range("O12").Value = "R$ 713.64"
$ws.Range("P12").Value = "R$ 713.64"

# Plain numeric cells.
$ws.Range("M6").Value = 30
$ws.Range("M7").Value = 252
$ws.Range("N7").Value = 2
$ws.Range("M9").Value = 60
$ws.Range("M10").Value = 252
